# Extend the phone-directory sheet (Лист1) with a duplicated owner-info
# block (columns E:I) and three new rows (6,7 were added, plus the
# existing row 5's contact info is now mirrored into E:I), while the
# old single "owner name" cell (old column E) is relocated to column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 : move the title out of column E into column H ---------------
$ws.Range("E1").ClearContents()
$ws.Range("A1").Value = 89859708676
$ws.Range("H1").Value = "Автомобили с бробегом"

# --- Row 2 ----------------------------------------------------------------
$ws.Range("A2").Value = 89050397320
$ws.Range("B2").Value = "ПАО ""Вымпел-Коммуникации"""
$ws.Range("C2").Value = "Республика Татарстан"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "ПАО ""Вымпел-Коммуникации"""
$ws.Range("F2").Value = "Республика Татарстан"
$ws.Range("G2").Value = "Хозяин"
$ws.Range("H2").Value = "Хозяин"
$ws.Range("I2").Value = "M"

# --- Row 3 : drop the old D3/E3 cells, move the name to H3 ---------------
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("A3").Value = 84993907516
$ws.Range("H3").Value = "Гидра - Авто"

# --- Row 4 ------------------------------------------------------------------
$ws.Range("A4").Value = 89128439000
$ws.Range("B4").Value = "ПАО ""Мобильные ТелеСистемы"""
$ws.Range("C4").Value = "Оренбургская обл."
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "ПАО ""Мобильные ТелеСистемы"""
$ws.Range("F4").Value = "Оренбургская обл."
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "GM-AVTO автомобили с пробегом"

# --- Row 5 ------------------------------------------------------------------
$ws.Range("A5").Value = 89241086744
$ws.Range("B5").Value = "ПАО ""МегаФон"""
$ws.Range("C5").Value = "Хабаровский край"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "ПАО ""МегаФон"""
$ws.Range("F5").Value = "Хабаровский край"
$ws.Range("G5").Value = "Михаил"
$ws.Range("H5").Value = "Михаил"
$ws.Range("I5").Value = "M"

# --- Row 6 (new) --------------------------------------------------------
$ws.Range("A6").Value = 89241086745
$ws.Range("B6").Value = "ПАО ""МегаФон"""
$ws.Range("C6").Value = "Хабаровский край"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "ПАО ""МегаФон"""
$ws.Range("F6").Value = "Хабаровский край"
$ws.Range("G6").Value = "Марина"
$ws.Range("H6").Value = "Марина"
$ws.Range("I6").Value = "Ж"

# --- Row 7 (new) --------------------------------------------------------
$ws.Range("A7").Value = 89241086746
$ws.Range("B7").Value = "ПАО ""МегаФон"""
$ws.Range("C7").Value = "Хабаровский край"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = "ПАО ""МегаФон"""
$ws.Range("F7").Value = "Хабаровский край"
$ws.Range("G7").Value = "Илья"
$ws.Range("H7").Value = "Илья"
$ws.Range("I7").Value = "M"

# --- Selection left where the author last was editing ---------------------
$ws.Range("B8").Select()
